# "First Functionality ETL version"
# Appends the next batch of ETL-imported questions (rows 273-301) to the
# `pregunta` sheet, each with its four answer options (columns E/G/I/K) and
# the boolean "is correct" flags (columns F/H/J/L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pregunta")

# Columns: A=pregunta_id, B=texto_pregunta, C=materia_id, D=examen_id,
#          E=texto_pregunta_1, F=Opcion_1, G=texto_pregunta_2, H=Opcion_2,
#          I=texto_pregunta_3, J=Opcion_3, K=texto_pregunta_4, L=Opcion_4
$newRows = @(
    @{ r = 273; vals = @(272, 'Giordano  Bruno  sostiene  que  lo  que  se  entiende  por  bello  depende  de  la  apreciación  particular  del  individuo. Incluso sus valoraciones estéticas pueden variar según los cambios de su estado anímico. De lo anterior se deduce que', 10, 2, 'la psicología demuestra la imposibilidad del arte y la inexistencia de la belleza.', $false, 'o bello carece de una significación única válida para todos los seres humanos', $true, 'la psicología puede definir la belleza a partir del análisis del individuo', $false, 'lo que se conoce como belleza es simplemente un estado psicológico', $false) },
    @{ r = 274; vals = @(273, 'La estética, en el sentido de discurso filosófico aplicado al arte y a sus relaciones con el bien y el mal, es decir, aplicado a su finalidad, se remonta a Platón. Pero gracia a Croce, tiende a desarrollarse como disciplina autónoma. De lo anterior podemos deducir que la estética como discurso filosófico ha evoluciona, porque ', 10, 3, 'el arte ha sido abordado desde diferentes corrientes filosóficas', $true, 'sus problemas han sido separados de otras problemáticas filosóficas', $false, 'el pensamiento filosófico hace referencia a los tratados estéticos', $false, 'la estética como finalidad ha buscado encontrar el sentido de la vida', $false) },
    @{ r = 275; vals = @(274, 'Para  Sartre,  la  libertad  humana  es  la  facultad  de  hacer  lo  que  uno  quiera  sin  ninguna  limitación.  La  libertad formalmente  no  tiene  ningún  contenido  determinante,  no  está  sometida  a  ninguna  necesidad  o  ley.  Sin  embargo, Sartre afirma que la libertad es una condena, porque ', 10, 4, 'la única forma de limitarla es haciendo una consideración formal de la norma', $false, 'exige del hombre madurez frente a la forma en la que elige su modo de actuar', $false, 'ella lleva la moralidad al establecimiento de los límites', $false, $true, $false) },
    @{ r = 276; vals = @(275, 'Karl Philipp Moritz escribió en 1785 que lo bello tiene una finalidad interna que consiste en provocar placer por su perfección. Esto quiere decir, que del mismo modo como a la técnica le corresponde una finalidad externa que es la utilidad, las bellas artes tienen como objetivo último el placer. Según lo anterior', 10, 5, 'una obra de arte es perfecta cuando se percibe su imagen real', $false, 'la contemplación de la belleza provoca un gozo profundo', $true, 'se halla gozo tanto en lo bello como en lo útil', $false, 'el arte viene dado por la imitación de la naturaleza que le sirve al hombre', $false) },
    @{ r = 277; vals = @(276, 'Según Hegel, el estudio de lo bello como objeto del arte es un estudio independiente del estudio de la naturaleza. La creación artística no puede estudiarse como una rama de la ciencia natural. En otras palabras, la belleza artística es independiente de la belleza en la naturaleza. De acuerdo con lo anterior, Hegel critica  ', 10, 6, 'la concepción según la cual el estudio del arte y el estudio de la naturaleza son independientes', $true, 'todo intento de separar la estética de la ciencia natural', $false, 'toda concepción que subordine el arte al estudio de la naturaleza', $false, 'la concepción de belleza en el arte como algo independiente de lo bello natural ', $false) },
    @{ r = 278; vals = @(277, 'Algunos afirman que en los siglos XVII y XVIII, se creía que el arte y la estética gozaban de una condición interna similar a la de la naturaleza, porque', 10, 7, 'el arte es imitación de la naturaleza y por lo tanto, se rige por leyes similares', $false, 'las leyes de la estética también se subordinan a un único principio que es el de la imitación', $false, 'el arte es obra humana y por lo tanto, se rige por los mismos principios que rigen la ciencia', $true, 'para poder plasmar la complejidad de la naturaleza, es necesario conocer sus leyes', $false) },
    @{ r = 279; vals = @(278, 'Para Descartes, la certeza en el conocimiento era fundamental, razón por la cual sometió al análisis de la razón no sólo a la filosofía, ala lógica, a la física y a la psicología, sino también al arte. Por lo tanto, se puede afirmar que para Descartes  ', 10, 8, 'el arte debía cumplir con las pruebas de la razón para demostrar su verdad', $false, 'era necesario un análisis del arte para ver si su contenido era verdadero y esencial ', $false, 'los principios del arte eran absolutamente racionales e inamovibles', $false, 'se debía realizar una crítica del arte para saber si era ciencia', $true) },
    @{ r = 280; vals = @(279, 'Schelling muestra cómo a través del arte se pone al descubierto el interior del filósofo, mostrando que en el arte y en las creaciones artísticas tienen su encuentro la naturaleza y el espíritu, la ley y la libertad, el cuerpo y el alma, y que el filósofo se descubre en el arte cuando', 10, 9, 'reconoce en lo bello la infinitud propia del alma', $true, 'su obra lo lleva a la liberación total y a la participación de lo que traspasa los límites', $false, 'se reconoce como creador de una obra que se constituye en símbolo del infinito', $false, 'hace descender lo infinito a la finitud a través del símbolo', $false) },
    @{ r = 281; vals = @(280, 'Según Aristóteles “la poesía es más filosófica que la historia”. En efecto, la historia narra lo que ha sucedido. En cambio, la poesía cuenta lo que siempre puede suceder; nos enseña a ver lo universal en el obrar humano y cómo la filosofía expresa las inquietudes fundamentales del hombre. Esto significa que', 10, 10, 'la filosofía debería guiarse por la poesía ', $false, 'la poesía guarda en sí una perspectiva filosófica, porque se interroga por lo universal ', $true, 'la historia sólo es un recuento de hechos sin interrogarse sobre el significado del obrar', $false, 'la historia debería hacerse más poética', $false) },
    @{ r = 282; vals = @(281, 'La  Geografía  tiene  gran  importanciaprincipalmente para', 10, 10, 'conocer las costumbres rurales', $false, 'la organización del espacio terrestre', $true, 'evitar el daño de la capa de ozono', $false, 'limitar el desarrollo de la industria', $false) },
    @{ r = 283; vals = @(282, 'Giordano  Bruno  sostiene  que  lo  que  se  entiende  por  bello  depende  de  la  apreciación  particular  del  individuo. Incluso sus valoraciones estéticas pueden variar según los cambios de su estado anímico. De lo anterior se deduce que', 10, 9, 'la psicología demuestra la imposibilidad del arte y la inexistencia de la belleza.', $false, 'o bello carece de una significación única válida para todos los seres humanos', $true, 'la psicología puede definir la belleza a partir del análisis del individuo', $false, 'lo que se conoce como belleza es simplemente un estado psicológico', $false) },
    @{ r = 284; vals = @(283, 'La estética, en el sentido de discurso filosófico aplicado al arte y a sus relaciones con el bien y el mal, es decir, aplicado a su finalidad, se remonta a Platón. Pero gracia a Croce, tiende a desarrollarse como disciplina autónoma. De lo anterior podemos deducir que la estética como discurso filosófico ha evoluciona, porque ', 10, 8, 'el arte ha sido abordado desde diferentes corrientes filosóficas', $true, 'sus problemas han sido separados de otras problemáticas filosóficas', $false, 'el pensamiento filosófico hace referencia a los tratados estéticos', $false, 'la estética como finalidad ha buscado encontrar el sentido de la vida', $false) },
    @{ r = 285; vals = @(284, 'Para  Sartre,  la  libertad  humana  es  la  facultad  de  hacer  lo  que  uno  quiera  sin  ninguna  limitación.  La  libertad formalmente  no  tiene  ningún  contenido  determinante,  no  está  sometida  a  ninguna  necesidad  o  ley.  Sin  embargo, Sartre afirma que la libertad es una condena, porque ', 10, 7, 'la única forma de limitarla es haciendo una consideración formal de la norma', $false, 'exige del hombre madurez frente a la forma en la que elige su modo de actuar', $false, 'ella lleva la moralidad al establecimiento de los límites', $false, $true, $false) },
    @{ r = 286; vals = @(285, 'Karl Philipp Moritz escribió en 1785 que lo bello tiene una finalidad interna que consiste en provocar placer por su perfección. Esto quiere decir, que del mismo modo como a la técnica le corresponde una finalidad externa que es la utilidad, las bellas artes tienen como objetivo último el placer. Según lo anterior', 10, 6, 'una obra de arte es perfecta cuando se percibe su imagen real', $false, 'la contemplación de la belleza provoca un gozo profundo', $true, 'se halla gozo tanto en lo bello como en lo útil', $false, 'el arte viene dado por la imitación de la naturaleza que le sirve al hombre', $false) },
    @{ r = 287; vals = @(286, 'Según Hegel, el estudio de lo bello como objeto del arte es un estudio independiente del estudio de la naturaleza. La creación artística no puede estudiarse como una rama de la ciencia natural. En otras palabras, la belleza artística es independiente de la belleza en la naturaleza. De acuerdo con lo anterior, Hegel critica  ', 10, 5, 'la concepción según la cual el estudio del arte y el estudio de la naturaleza son independientes', $true, 'todo intento de separar la estética de la ciencia natural', $false, 'toda concepción que subordine el arte al estudio de la naturaleza', $false, 'la concepción de belleza en el arte como algo independiente de lo bello natural ', $false) },
    @{ r = 288; vals = @(287, 'Algunos afirman que en los siglos XVII y XVIII, se creía que el arte y la estética gozaban de una condición interna similar a la de la naturaleza, porque', 10, 4, 'el arte es imitación de la naturaleza y por lo tanto, se rige por leyes similares', $false, 'las leyes de la estética también se subordinan a un único principio que es el de la imitación', $false, 'el arte es obra humana y por lo tanto, se rige por los mismos principios que rigen la ciencia', $true, 'para poder plasmar la complejidad de la naturaleza, es necesario conocer sus leyes', $false) },
    @{ r = 289; vals = @(288, 'Para Descartes, la certeza en el conocimiento era fundamental, razón por la cual sometió al análisis de la razón no sólo a la filosofía, ala lógica, a la física y a la psicología, sino también al arte. Por lo tanto, se puede afirmar que para Descartes  ', 10, 3, 'el arte debía cumplir con las pruebas de la razón para demostrar su verdad', $false, 'era necesario un análisis del arte para ver si su contenido era verdadero y esencial ', $false, 'los principios del arte eran absolutamente racionales e inamovibles', $false, 'se debía realizar una crítica del arte para saber si era ciencia', $true) },
    @{ r = 290; vals = @(289, 'Schelling muestra cómo a través del arte se pone al descubierto el interior del filósofo, mostrando que en el arte y en las creaciones artísticas tienen su encuentro la naturaleza y el espíritu, la ley y la libertad, el cuerpo y el alma, y que el filósofo se descubre en el arte cuando', 10, 2, 'reconoce en lo bello la infinitud propia del alma', $true, 'su obra lo lleva a la liberación total y a la participación de lo que traspasa los límites', $false, 'se reconoce como creador de una obra que se constituye en símbolo del infinito', $false, 'hace descender lo infinito a la finitud a través del símbolo', $false) },
    @{ r = 291; vals = @(290, 'Según Aristóteles “la poesía es más filosófica que la historia”. En efecto, la historia narra lo que ha sucedido. En cambio, la poesía cuenta lo que siempre puede suceder; nos enseña a ver lo universal en el obrar humano y cómo la filosofía expresa las inquietudes fundamentales del hombre. Esto significa que', 10, 1, 'la filosofía debería guiarse por la poesía ', $false, 'la poesía guarda en sí una perspectiva filosófica, porque se interroga por lo universal ', $true, 'la historia sólo es un recuento de hechos sin interrogarse sobre el significado del obrar', $false, 'la historia debería hacerse más poética', $false) },
    @{ r = 292; vals = @(291, 'La Geografía Regional, es importante para', 9, 10, 'entender  la  universalidad  de  losfenómenos geográficosfenómenos geográfico', $false, 'visualizar globalmente globalmente el el espacioespacioterrestreterrestr', $false, 'concatenar   todo   los   fenómenosplanetariosplanetario', $false, 'identificar  las  potencialidades  delespacio terrestreespacio terrestr', $true) },
    @{ r = 293; vals = @(292, 'Las  dos  grandes  doctrinas  geográficas  queexplican  el  desarrollo  cultural  y  económicoexplican  el  desarrollo  cultural  y  económicode los pueblos, son:de los pueblos, son', 9, 9, 'el idealismo y el materialismo', $false, 'el positivismo y el historicismo', $false, 'el determinismo y el posibilismo', $true, 'el anarquismo y el conservadurismo', $false) },
    @{ r = 294; vals = @(293, 'La deforestación de las vertientesaltoandinas   está   provocando   huaycos   ealtoandinas   está   provocando   huaycos   einundaciones    en    la    selva    alta.    ¿Quéinundaciones    en    la    selva    alta.    ¿Quéprincipio geográfico se aplica?principio geográfico se aplica?', 9, 8, 'Localización', $false, 'Descripción', $false, 'Causalidad', $true, 'Comparación', $false) },
    @{ r = 295; vals = @(294, 'La edafología, es una rama de la GeografíaFísica que estudia:Física que estudia', 9, 7, 'la distribución espacial del agua', $false, 'el origen y la evolución de los suelos', $true, 'las variadas morfologías del terreno', $false, 'el   grado   de   vulnerabilidad   de   losecosistemas.', $false) },
    @{ r = 296; vals = @(295, '¿De qué depende que el impacto ambientaldel  hombre  sobre  el  espacio  geográfico  seadel  hombre  sobre  el  espacio  geográfico  seanegativo o positivo?negativo o positivo?', 9, 6, 'Depende esencialmente de sucondición económica y culturalcondición económica y cultura', $true, 'Depende  de  su  grado  de  concienciaecológicaecológica', $false, 'Depende de su estado de animo', $false, 'Depende     de     sus     característicasgenéticas y antropológicasgenéticas y antropológi', $false) },
    @{ r = 297; vals = @(296, 'La geografía regional tiene como objeto deestudio', 9, 5, 'los fenómenos de alcance planetario', $false, 'una zona espacial bien definida', $true, 'la contaminación ambiental', $false, 'la  universalidad  de  los  fenómenosterrestres', $false) },
    @{ r = 298; vals = @(297, 'Con  los  aportes  de  Alexander  vonHumboldt  y  Paúl  Vidal  de  la  Blache  lageografía  se  convierte  fundamentalmenteen una ciencia', 9, 4, 'que describe los paisajes regionales', $false, 'que   cuantifica   las   dimensionesespaciales', $false, 'holística  que  explica  la  relaciónhombre-naturaleza', $true, 'que rinde pleitesía a la naturaleza', $false) },
    @{ r = 299; vals = @(298, 'La  causalidad  es  un  principio  dialécticoaplicado a la geografía por', 9, 3, 'Paúl Vidal de la Blache', $false, 'Jean Brunhes', $false, 'Alexander von Humboldt', $true, 'Federico Ratzell', $false) },
    @{ r = 300; vals = @(299, 'Para  la  ciencia  geográfica  el  hombre  esprincipalmente un', 9, 2, 'depredador innato de la naturaleza', $false, 'agente  transformador  del  mediofísico', $true, 'apasionado de la naturaleza', $false, 'activista    de    los    problemasambientales', $false) },
    @{ r = 301; vals = @(300, 'En cuál de los siguientes casos se aplica elprincipio de comparación', 9, 1, 'La  costa  peruana  es  una  planicieangosta y desértica', $false, 'La  costa  peruana  se  extiende  desdeTacna a Tumbes', $false, 'La costa peruana a diferencia la deEcuador es desértica', $true, 'La costa peruana es desértico debidoal Mar Frío', $false) }
)

foreach ($item in $newRows) {
    $rowNum = $item.r
    $vals = $item.vals
    for ($col = 1; $col -le $vals.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value2 = $vals[$col - 1]
    }
}

# Column K (4th answer option) uses a distinct text style whenever it holds
# option text rather than a bare boolean (mirrors the existing rows above).
# Copy that formatting from an already-styled K cell so no new style is
# introduced.
$ws.Range("K262").Copy()
foreach ($item in $newRows) {
    $kVal = $item.vals[10]
    if ($kVal -is [string]) {
        $ws.Cells.Item($item.r, 11).PasteSpecial(-4122)
    }
}
$excel.CutCopyMode = $false

# Refresh the sheet view to the new zoom level / scroll position.
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
$ws.Range("D1").Select()

Write-Host "Added rows 273-301 ($($newRows.Count) questions) to 'pregunta' sheet"
